$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 381, shifting existing rows (381-433) down to (382-434)
$ws.Rows.Item(381).Insert()

# Populate the newly inserted row 381 with the new weekly record
$ws.Cells.Item(381, 1).Value = 4
$ws.Cells.Item(381, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(381, 3).Value = "Los Lagos"
$ws.Cells.Item(381, 4).Value = 45127
$ws.Cells.Item(381, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(381, 5).Value = 10
$ws.Cells.Item(381, 6).Value = 100112032
$ws.Cells.Item(381, 7).Value = "Zapallo italiano"
$ws.Cells.Item(381, 8).Value = "Sin especificar"
$ws.Cells.Item(381, 9).Value = "Primera"
$ws.Cells.Item(381, 10).Value = 120
$ws.Cells.Item(381, 11).Value = 19000
$ws.Cells.Item(381, 12).Value = 19000
$ws.Cells.Item(381, 13).Value = 19000
$ws.Cells.Item(381, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(381, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(381, 16).Value = 380
$ws.Cells.Item(381, 17).Value = 50
$ws.Cells.Item(381, 18).Value = "Hortaliza"
